$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'89.723.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.31%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.178.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.35%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("E5").Value = "'  -0.89%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'614.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.25%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.389"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.69%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.690"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.93%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'3.176.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.29%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.576"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.63%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  -5.38%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  -7.96%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'3.770.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.16%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'89.536.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.44%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'32.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -6.12%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'5.25"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.209.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.11%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'3.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.46%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'13.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -6.23%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'435.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.67%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.0000195"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +33.95%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'8.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.53%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'5.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -6.60%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'5.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.23%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'11.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -6.43%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'3.349.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.91%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'75.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.62%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  -8.48%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.36%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'4.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +26.01%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'8.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.56%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'533.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.80%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'7.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.10%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  -6.86%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  -10.02%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "'EthereumClassic"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'21.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.50%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("B39").Value = "'WhiteBITCoin"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'22.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.12%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  -9.83%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.12%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  -0.05%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'1.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -7.17%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.371"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -9.09%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'148.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.18%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'  -3.44%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'172.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.43%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.124"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -10.08%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -9.67%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'4.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.87%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.608"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.06%  "
$ws.Range("E51").Style = "Normal"
